$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.111.43'
$ws.Range("D3").Value = '3.524.90'
$ws.Range("E3").Value = '  +2.38%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '''596.55'
$ws.Range("E5").Value = '  +1.18%  '
$ws.Range("D6").Value = '''137.92'
$ws.Range("E6").Value = '  -0.24%  '
$ws.Range("D7").Value = '3.525.57'
$ws.Range("E7").Value = '  +2.44%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("E9").Value = '  -1.46%  '
$ws.Range("E10").Value = '  +1.78%  '
$ws.Range("D11").Value = '''6.80'
$ws.Range("E11").Value = '  -7.40%  '
$ws.Range("D12").Value = '''0.387'
$ws.Range("E12").Value = '  +2.21%  '
$ws.Range("D13").Value = '4.125.00'
$ws.Range("E13").Value = '  +2.45%  '
$ws.Range("E14").Value = '  +1.39%  '
$ws.Range("E15").Value = '  +1.89%  '
$ws.Range("D16").Value = '3.522.13'
$ws.Range("E16").Value = '  +2.19%  '
$ws.Range("E17").Value = '  +1.39%  '
$ws.Range("D18").Value = '65.190.51'
$ws.Range("E18").Value = '  -0.59%  '
$ws.Range("D19").Value = '''10.25'
$ws.Range("E19").Value = '  +3.47%  '
$ws.Range("D20").Value = '''5.95'
$ws.Range("E20").Value = '  +0.91%  '
$ws.Range("D21").Value = '''14.24'
$ws.Range("E21").Value = '  +3.65%  '
$ws.Range("D22").Value = '''391.63'
$ws.Range("E22").Value = '  -0.43%  '
$ws.Range("D23").Value = '''0.570'
$ws.Range("E23").Value = '  +2.43%  '
$ws.Range("D24").Value = '3.666.39'
$ws.Range("E24").Value = '  +2.46%  '
$ws.Range("D25").Value = '''73.73'
$ws.Range("E25").Value = '  +0.40%  '
$ws.Range("D26").Value = '''1.00'
$ws.Range("E26").Value = '  -0.01%  '
$ws.Range("E27").Value = '  +6.23%  '
$ws.Range("E28").Value = '  +5.95%  '
$ws.Range("D29").Value = '''0.999'
$ws.Range("E29").Value = '  -0.31%  '
$ws.Range("E30").Value = '  +2.09%  '
$ws.Range("E31").Value = '  -1.62%  '
$ws.Range("D32").Value = '3.539.99'
$ws.Range("E32").Value = '  +2.66%  '
$ws.Range("E33").Value = '  -0.01%  '
$ws.Range("D34").Value = '''23.74'
$ws.Range("E34").Value = '  +3.01%  '
$ws.Range("D35").Value = '''0.144'
$ws.Range("E35").Value = '  -1.26%  '
$ws.Range("E36").Value = '  +6.63%  '
$ws.Range("E37").Value = '  +0.11%  '
$ws.Range("D38").Value = '''168.64'
$ws.Range("E38").Value = '  -2.54%  '
$ws.Range("E39").Value = '  +4.04%  '
$ws.Range("D40").Value = '''4.95'
$ws.Range("E40").Value = '  +2.77%  '
$ws.Range("E41").Value = '  +3.98%  '
$ws.Range("E42").Value = '  -0.57%  '
$ws.Range("D43").Value = '''25.83'
$ws.Range("E43").Value = '  +13.42%  '
$ws.Range("D44").Value = '''42.74'
$ws.Range("E44").Value = '  -2.28%  '
$ws.Range("E45").Value = '  -0.02%  '
$ws.Range("D46").Value = '''4.40'
$ws.Range("E46").Value = '  -0.79%  '
$ws.Range("E47").Value = '  +2.32%  '
$ws.Range("E48").Value = '  +4.36%  '
$ws.Range("D49").Value = '''6.78'
$ws.Range("E49").Value = '  +3.16%  '
$ws.Range("D50").Value = '2.372.92'
$ws.Range("E50").Value = '  +7.14%  '
$ws.Range("D51").Value = '''300.92'
